# Applies the edits described in the commit diff to Jacob Script.docx
$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2

# 1) Intro paragraph rewrite
$old1 = "Hello. My name is Jacob Dyer, and this is a portion of my portion of the final video for CS 560. About halfway Luke Duball will take over and then I will return to wrap up our discussion. "
$new1 = "Hello. My name is Jacob Dyer, and in this video Lukas Duball and myself are going to talk about our final project for CS 524. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) Insert "that exist in this universe" before the period ending "...within the factions."
$old2 = "to become powerful within the factions. These factions"
$new2 = "to become powerful within the factions that exist in this universe. These factions"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) Change "42 gigs" to "3 gigs" and add "around those bodies" before the trailing period
$old3 = "I think I read it was last about 42 gigs. The files we used cover less about the natural elements, such as celestial bodies, and more about the civilization that has been made. "
$new3 = "I think I read it was last about 3 gigs. The files we used cover less about the natural elements, such as celestial bodies, and more about the civilization that has been made around those bodies. "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null
